# Update cryptos list with fresh price/volume figures (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.379.60"
$ws.Range("E2").Value = "  -0.04%  "

$ws.Range("D3").Value = "2.068.23"
$ws.Range("E3").Value = "  +0.20%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.33"
$ws.Range("E5").Value = "  -0.25%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.624"
$ws.Range("E6").Value = "  +1.67%  "

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "57.18"
$ws.Range("E8").Value = "  -1.87%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.394"
$ws.Range("E9").Value = "  +2.82%  "

$ws.Range("E10").Value = "  +1.18%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.102"
$ws.Range("E11").Value = "  +0.58%  "

$ws.Range("D12").Value = "2.372.31"
$ws.Range("E12").Value = "  +0.23%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.37"
$ws.Range("E13").Value = "  -1.55%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.65"
$ws.Range("E14").Value = "  -1.45%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.775"
$ws.Range("E15").Value = "  -0.47%  "

$ws.Range("E16").Value = "  -0.49%  "

$ws.Range("D17").Value = "2.068.75"
$ws.Range("E17").Value = "  +0.28%  "

$ws.Range("D18").Value = "37.305.22"
$ws.Range("E18").Value = "  -0.77%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.17"
$ws.Range("E19").Value = "  -0.68%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.52"
$ws.Range("E20").Value = "  +0.64%  "

$ws.Range("E21").Value = "  -0.15%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "226.67"
$ws.Range("E22").Value = "  -0.01%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.43"
$ws.Range("E24").Value = "  +1.98%  "

$ws.Range("E25").Value = "  -2.55%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "166.85"
$ws.Range("E26").Value = "  +1.32%  "

$ws.Range("E27").Value = "  +0.99%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.41"
$ws.Range("E28").Value = "  -5.65%  "

$ws.Range("E29").Value = "  +0.99%  "

$ws.Range("E30").Value = "  -0.79%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.117"
$ws.Range("E31").Value = "  -1.79%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.52"
$ws.Range("E32").Value = "  +0.56%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0616"
$ws.Range("E33").Value = "  -1.33%  "

$ws.Range("E34").Value = "  +0.55%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.47"
$ws.Range("E35").Value = "  -3.21%  "

$ws.Range("E36").Value = "  +0.20%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.32"
$ws.Range("E37").Value = "  -3.01%  "

$ws.Range("E38").Value = "  -0.02%  "

$ws.Range("E39").Value = "  -4.84%  "

$ws.Range("E40").Value = "  -0.66%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0955"
$ws.Range("E41").Value = "  -3.14%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "97.65"

$ws.Range("D43").Value = "1.481.77"
$ws.Range("E43").Value = "  +0.38%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0212"
$ws.Range("E44").Value = "  +0.88%  "

$ws.Range("E45").Value = "  -0.48%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.17"
$ws.Range("E46").Value = "  -7.70%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.02"
$ws.Range("E47").Value = "  -0.39%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.20"
$ws.Range("E48").Value = "  -1.44%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "15.05"
$ws.Range("E49").Value = "  -5.68%  "

$ws.Range("E50").Value = "  +1.02%  "

$ws.Range("B51").Value = "MultiversX"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "47.88"
$ws.Range("E51").Value = "  +6.89%  "

Write-Output "Applied cryptos update"
